# Refresh the cryptocurrency price table (rows 2-51) with the latest scrape.
# Each row entry below lists only the columns that changed for that coin row:
#   B = Coin name, C = Link, D = Price, E = Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D='29.570.12'; E='  -1.44%  ' }
    @{ Row=3; D='2.023.93'; E='  -4.49%  ' }
    @{ Row=4; D='1.020'; E='  +1.01%  ' }
    @{ Row=5; D='332.13'; E='  -4.39%  ' }
    @{ Row=6; D='1.014'; E='  +0.49%  ' }
    @{ Row=7; D='0.4938'; E='  -5.14%  ' }
    @{ Row=8; D='0.4162'; E='  -6.65%  ' }
    @{ Row=9; D='53.80'; E='  -0.33%  ' }
    @{ Row=10; D='0.08795'; E='  -5.90%  ' }
    @{ Row=11; D='1.118'; E='  -5.51%  ' }
    @{ Row=12; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='2.076.68'; E='  -2.93%  ' }
    @{ Row=13; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='22.63'; E='  -10.29%  ' }
    @{ Row=14; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='8.138'; E='  -4.77%  ' }
    @{ Row=15; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='6.489'; E='  -6.12%  ' }
    @{ Row=16; D='97.02'; E='  -5.71%  ' }
    @{ Row=17; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='1.024'; E='  +1.33%  ' }
    @{ Row=18; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.00001102'; E='  -5.24%  ' }
    @{ Row=19; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.06637'; E='  -1.00%  ' }
    @{ Row=20; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='19.26'; E='  -10.86%  ' }
    @{ Row=21; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.008'; E='  -0.07%  ' }
    @{ Row=22; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='5.959'; E='  -5.64%  ' }
    @{ Row=23; D='29.657.04'; E='  -1.31%  ' }
    @{ Row=24; D='11.76'; E='  -7.66%  ' }
    @{ Row=25; D='2.319'; E='  -0.39%  ' }
    @{ Row=26; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='2.264.20'; E='  -5.39%  ' }
    @{ Row=27; D='158.64'; E='  -2.62%  ' }
    @{ Row=28; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='20.53'; E='  -7.26%  ' }
    @{ Row=29; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='6.438'; E='  -1.63%  ' }
    @{ Row=30; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='2.301'; E='  -9.59%  ' }
    @{ Row=31; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='127.48'; E='  -5.06%  ' }
    @{ Row=32; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='1.037'; E='  -10.14%  ' }
    @{ Row=33; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.09768'; E='  -7.57%  ' }
    @{ Row=34; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.529'; E='  -14.15%  ' }
    @{ Row=35; D='3.870'; E='  -2.73%  ' }
    @{ Row=36; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='5.786'; E='  -7.54%  ' }
    @{ Row=37; D='9.676'; E='  -10.07%  ' }
    @{ Row=38; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.02424'; E='  -7.48%  ' }
    @{ Row=39; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='1.303'; E='  -1.95%  ' }
    @{ Row=40; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.06293'; E='  -8.65%  ' }
    @{ Row=41; D='11.74'; E='  -7.65%  ' }
    @{ Row=42; B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.6418'; E='  -9.32%  ' }
    @{ Row=43; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.2040'; E='  -9.03%  ' }
    @{ Row=44; B='Frax'; C='https://coinranking.com/coin/KfWtaeV1W+frax-frax'; D='1.011'; E='  +0.25%  ' }
    @{ Row=45; B='Decentraland'; C='https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; D='0.6283'; E='  -8.35%  ' }
    @{ Row=46; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='13.70'; E='  -5.87%  ' }
    @{ Row=47; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='2.185'; E='  -7.91%  ' }
    @{ Row=48; B='WEMIXTOKEN'; C='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; D='1.279'; E='  +0.29%  ' }
    @{ Row=49; B='PancakeSwap'; C='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D='3.586'; E='  -1.48%  ' }
    @{ Row=50; B='BabyDogeCoin'; C='https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; D='0.00000000334'; E='  -8.31%  ' }
    @{ Row=51; D='0.07067'; E='  -0.56%  ' }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
    if ($u.ContainsKey("D")) {
        $cell = $ws.Cells.Item($u.Row, 4)
        # Prices such as "1.020" or "53.80" parse as numbers and would lose their
        # literal form (trailing zeros) unless the cell is forced to Text first.
        if ($u.D -match "^-?\d+(\.\d+)?$") {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
}
